# [Refactor] Added additional failing files to test
#
# Writes the "some data" header into A1 of every worksheet (setPart1,
# setPart2, setPart3) with an explicit font so a second font/cellXf pair is
# minted, fixes up sheet3's view (it was missing <sheetViews> and carried a
# stray " NA " text node), and moves the active tab to the third sheet.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("A1")
    $cell.Value = "some data"
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
}

# Make the third sheet (setPart3) the active tab (activeTab index 2).
$wb.Worksheets.Item(3).Activate()
